$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 433
$ws.Range("B2").Value = "SCC"
$ws.Range("C2").Value = "Property & Construction"
$ws.Range("D2").Value = "SET50 / SETCLMV / SETHD / SETTHSI"
